$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date for rows 2-11 from 2023-10-22
# to 2023-10-25, keeping the existing date cell formatting.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = "2023-10-25"
}
